$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the sheets
$ws1.Name = "pEqui"
$ws2.Name = "interestRates"

# --- interestRates (sheet2): replace contents with the new interest-rate table ---
$ws2.Range("A1").Value = "Rrate/t"
$ws2.Range("B1").Value = "Rrate/Rrate"

$ws2.Range("A2").Value = 2016
$ws2.Range("B2").Value = 0.94788653507940634

$ws2.Range("A3").Value = 2017
$ws2.Range("B3").Value = 1.0248619914586095

$ws2.Range("A4").Value = 2018
$ws2.Range("B4").Value = 1.1545674338403451

$ws2.Range("A5").Value = 2019
$ws2.Range("B5").Value = 1.0375000000000001

$ws2.Range("A6").Value = 2020
$ws2.Range("B6").Value = 1.0409999999999997

$ws2.Range("A7").Value = 2021
$ws2.Range("B7").Value = 1.0450000000000002

$ws2.Range("A8").Value = 2022
$ws2.Range("B8").Value = 1.0480000000000003

$ws2.Range("A9").Value = 2023
$ws2.Range("B9").Value = 1.05

$ws2.Range("A10").Value = 2024
$ws2.Range("B10").Value = 1.0520000000000003

$ws2.Range("A11").Value = 2025
$ws2.Range("B11").Value = 1.054

$ws2.Range("A12").Value = 2026
$ws2.Range("B12").Value = 1.0560000000000003

$ws2.Range("A13").Value = 2027
$ws2.Range("B13").Value = 1.0580000000000001

$ws2.Range("A14").Value = 2028
$ws2.Range("B14").Value = 1.0589999999999999

$ws2.Range("A15").Value = 2029
$ws2.Range("B15").Value = 1.0609999999999999

$ws2.Range("A16").Value = 2030
$ws2.Range("B16").Value = 1.0619999999999998

# interestRates page setup (paperSize/orientation match pEqui's)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Set the selection on interestRates (this also makes it the active sheet momentarily)
$ws2.Range("U15").Select()

# --- pEqui (sheet1): append two new rows ---
$ws1.Range("A7").Value = "I_iB"
$ws1.Range("B7").Value = 1

$ws1.Range("A8").Value = "I_iM"
$ws1.Range("B8").Value = 1

# Re-select pEqui as the active sheet/cell (matches final tabSelected + selection state)
$ws1.Range("A2").Select()
